$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'295.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.61%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.72%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.044"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.31%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07487"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.23%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.357"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.80%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.45%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9266"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'0.1192"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.51%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1819"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.90%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08904"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.81%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04184"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.09%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.03%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001278"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.20%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005893"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.39%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.346"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.58%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3311"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.51%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.907"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.11%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'4.09%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2968"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.66%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04045"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.13%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003862"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.88%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001231"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-4.18%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.34%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02393"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'5.64%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05203"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'5.75%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006587"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.68%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.04%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'4.45%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007381"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.34%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008118"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'16.95%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3221"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.37%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006209"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.76%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.30%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'0.004203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.30%"
$ws.Range("E51").Style = "Normal"
